$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = -0.2809526231947743
$ws.Range("J18").Value = 0.1328999324788512
$ws.Range("K18").Value = -0.01258305283788508
$ws.Range("L18").Value = 1.984619266005203
